$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C39").Value = "https://www.amazon.de/dp/B078XXP3RH?ref=cm_sw_r_cso_cp_apin_dp_CR2DCDFB6A4A5F0S67N2&social_share=cm_sw_r_cso_cp_apin_dp_CR2DCDFB6A4A5F0S67N2&badgeInsights=bestseller-insights&th=1"
$ws.Range("B39").Value = "https://m.media-amazon.com/images/I/81NlaW85aBL._AC_SL1500_.jpg"
$ws.Range("A39").Value = "Gantere 2kg"
$ws.Range("D39").Value = "18 EUR"

$ws.Range("D40").Select()
